$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 77, shifting existing rows 77-124 down to 78-125.
$ws.Rows.Item(77).Insert()

# Populate the newly inserted row 77 with its data.
$ws.Range("A77").Value = 3
$ws.Range("B77").Value = "Femacal de La Calera"
$ws.Range("C77").Value = "Coquimbo"
$ws.Range("D77").Value = 44582
$ws.Range("E77").Value = 5
$ws.Range("F77").Value = 100112052
$ws.Range("G77").Value = "Albahaca"
$ws.Range("H77").Value = "Sin especificar"
$ws.Range("I77").Value = "Primera"
$ws.Range("J77").Value = 130
$ws.Range("K77").Value = 4000
$ws.Range("L77").Value = 4500
$ws.Range("M77").Value = 4269
$ws.Range("N77").Value = "$/docena de matas"
$ws.Range("O77").Value = "Provincia de Quillota"
$ws.Range("P77").Value = 712
$ws.Range("Q77").Value = 6
$ws.Range("R77").Value = "Hortaliza"

Write-Output "Row inserted and populated"
